# Update the two "smart folder" test rows (anulación / obtener número de anulación)
# in the single sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'0420194406824"
$ws.Range("C2").Value = "'32581"

$ws.Range("B3").Value = "'0420172008636  "
$ws.Range("C3").Value = "'55299   "

$ws.Range("C7").Select() | Out-Null
